# Apply the row permutation described by the diff.
# The data rows (2-5) get their values rotated:
#   old Row 2 -> new Row 4
#   old Row 3 -> new Row 5
#   old Row 4 -> new Row 3
#   old Row 5 -> new Row 2
# Columns A,B,C,E,F,G,H,O,R are identical across all rows, so only the
# columns that actually vary (D,I,J,K,L,M,N,P,Q) need to be written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values for the columns that differ row to row,
# before we start overwriting anything.
$cols = @("D","I","J","K","L","M","N","P","Q")

$orig = @{}
foreach ($r in 2..5) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $rowVals
}

# Mapping: which original row's data should end up in each new row.
$mapping = @{ 2 = 5; 3 = 4; 4 = 2; 5 = 3 }

foreach ($newRow in 2..5) {
    $srcRow = $mapping[$newRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $orig[$srcRow][$c]
    }
}

$wb.Save()
